# Regenerate save_data to use K (column G) instead of Strike#.
# Recalculated std/mean, and write the new K (s_vals) values into column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(4, 0, 1, 1, 3, 1, 2, 0, 1, 3, 0, 2, 0, 2, 2, 3, 2, 1, 0)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
